$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")
$ws.Range("N1").Value = "description"
$null = $ws.Range("N1").Select()
